$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format date-text cells (Y/AA columns store dates as literal text,
# e.g. "2023-09-06"). Excel's COM layer auto-parses ISO-date-looking
# strings into date serials on assignment, so force Text format first
# to preserve the original literal-string semantics.
$dateTextCells = @(
    "Y46",
    "AA46",
    "Y47",
    "AA47",
    "Y48",
    "AA48",
    "Y52",
    "AA52",
    "Y65",
    "AA65",
    "Y66",
    "AA66",
    "Y67",
    "AA67",
    "Y68",
    "AA68"
)
foreach ($addr in $dateTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

    # Row 43
    $ws.Range("A43").Value = 112080974
    $ws.Range("Q43").Value = 838233
    $ws.Range("R43").Value = 7446396

    # Row 44
    $ws.Range("A44").Value = 112080972
    $ws.Range("Q44").Value = 838237
    $ws.Range("R44").Value = 7446403

    # Row 45
    $ws.Range("A45").Value = 112080973
    $ws.Range("B45").Value = 86357
    $ws.Range("E45").Value = 4412
    $ws.Range("F45").Value = 'Äggvaxskivling'
    $ws.Range("G45").Value = 'Hygrophorus karstenii'
    $ws.Range("H45").Value = 'Sacc. & Cub.'
    $ws.Range("Q45").Value = 838236
    $ws.Range("R45").Value = 7446398
    $ws.Range("AF45").Value = ""

    # Row 46
    $ws.Range("A46").Value = 112080969
    $ws.Range("B46").Value = 86357
    $ws.Range("E46").Value = 4412
    $ws.Range("F46").Value = 'Äggvaxskivling'
    $ws.Range("G46").Value = 'Hygrophorus karstenii'
    $ws.Range("H46").Value = 'Sacc. & Cub.'
    $ws.Range("Q46").Value = 838168
    $ws.Range("R46").Value = 7446363
    $ws.Range("Y46").Value = '2023-09-06'
    $ws.Range("AA46").Value = '2023-09-06'

    # Row 47
    $ws.Range("A47").Value = 112081019
    $ws.Range("B47").Value = 78633
    $ws.Range("D47").Value = 'LC'
    $ws.Range("E47").Value = 6456
    $ws.Range("F47").Value = 'Skinnlav'
    $ws.Range("G47").Value = 'Leptogium saturninum'
    $ws.Range("H47").Value = '(Dicks.) Nyl.'
    $ws.Range("Q47").Value = 838202
    $ws.Range("R47").Value = 7446426
    $ws.Range("Y47").Value = '2023-09-06'
    $ws.Range("AA47").Value = '2023-09-06'
    $ws.Range("AF47").Value = ""

    # Row 48
    $ws.Range("A48").Value = 112081047
    $ws.Range("B48").Value = 90792
    $ws.Range("E48").Value = 4361
    $ws.Range("F48").Value = 'Orange taggsvamp'
    $ws.Range("G48").Value = 'Hydnellum aurantiacum'
    $ws.Range("H48").Value = '(Batsch:Fr.) P.Karst.'
    $ws.Range("Q48").Value = 837882
    $ws.Range("R48").Value = 7446387
    $ws.Range("Y48").Value = '2023-08-31'
    $ws.Range("AA48").Value = '2023-08-31'

    # Row 49
    $ws.Range("A49").Value = 112081061
    $ws.Range("B49").Value = 88166
    $ws.Range("D49").Value = 'VU'
    $ws.Range("E49").Value = 6276
    $ws.Range("F49").Value = 'Goliatmusseron'
    $ws.Range("G49").Value = 'Tricholoma matsutake'
    $ws.Range("H49").Value = '(S.Ito & S.Imai) Singer'
    $ws.Range("Q49").Value = 837896
    $ws.Range("R49").Value = 7446211
    $ws.Range("AF49").Value = ''

    # Row 50
    $ws.Range("A50").Value = 112081070
    $ws.Range("B50").Value = 81371
    $ws.Range("D50").Value = 'NT'
    $ws.Range("E50").Value = 1312
    $ws.Range("F50").Value = 'Gammelgransskål'
    $ws.Range("G50").Value = 'Pseudographis pinicola'
    $ws.Range("H50").Value = '(Nyl.) Rehm'
    $ws.Range("Q50").Value = 838037
    $ws.Range("R50").Value = 7446422
    $ws.Range("AF50").Value = ''

    # Row 51
    $ws.Range("A51").Value = 112080991
    $ws.Range("B51").Value = 56575
    $ws.Range("E51").Value = 103021
    $ws.Range("F51").Value = 'Talltita'
    $ws.Range("G51").Value = 'Poecile montanus'
    $ws.Range("H51").Value = '(Conrad von Baldenstein, 1827)'
    $ws.Range("Q51").Value = 837924
    $ws.Range("R51").Value = 7446271

    # Row 52
    $ws.Range("A52").Value = 112080990
    $ws.Range("B52").Value = 89559
    $ws.Range("D52").Value = 'NT'
    $ws.Range("E52").Value = 5442
    $ws.Range("F52").Value = 'Tallticka'
    $ws.Range("G52").Value = 'Porodaedalea pini'
    $ws.Range("H52").Value = '(Brot.) Murrill'
    $ws.Range("Q52").Value = 837871
    $ws.Range("R52").Value = 7446492
    $ws.Range("Y52").Value = '2023-08-31'
    $ws.Range("AA52").Value = '2023-08-31'

    # Row 60
    $ws.Range("A60").Value = 112081043
    $ws.Range("B60").Value = 95693
    $ws.Range("E60").Value = 221941
    $ws.Range("F60").Value = 'Plattlummer'
    $ws.Range("G60").Value = 'Lycopodium complanatum'
    $ws.Range("H60").Value = 'L.'
    $ws.Range("Q60").Value = 837847
    $ws.Range("R60").Value = 7446312

    # Row 61
    $ws.Range("A61").Value = 112081045
    $ws.Range("B61").Value = 90792
    $ws.Range("D61").Value = 'NT'
    $ws.Range("E61").Value = 4361
    $ws.Range("F61").Value = 'Orange taggsvamp'
    $ws.Range("G61").Value = 'Hydnellum aurantiacum'
    $ws.Range("H61").Value = '(Batsch:Fr.) P.Karst.'
    $ws.Range("Q61").Value = 837916
    $ws.Range("R61").Value = 7446236

    # Row 62
    $ws.Range("A62").Value = 112081028
    $ws.Range("B62").Value = 90812
    $ws.Range("E62").Value = 4366
    $ws.Range("F62").Value = 'Skarp dropptaggsvamp'
    $ws.Range("G62").Value = 'Hydnellum peckii'
    $ws.Range("H62").Value = 'Banker'
    $ws.Range("Q62").Value = 837840
    $ws.Range("R62").Value = 7446337

    # Row 63
    $ws.Range("A63").Value = 112081012
    $ws.Range("B63").Value = 78726
    $ws.Range("D63").Value = 'LC'
    $ws.Range("E63").Value = 6462
    $ws.Range("F63").Value = 'Stuplav'
    $ws.Range("G63").Value = 'Nephroma bellum'
    $ws.Range("H63").Value = '(Spreng.) Tuck.'
    $ws.Range("Q63").Value = 838034
    $ws.Range("R63").Value = 7446280

    # Row 64
    $ws.Range("A64").Value = 112081020
    $ws.Range("B64").Value = 78633
    $ws.Range("D64").Value = 'LC'
    $ws.Range("E64").Value = 6456
    $ws.Range("F64").Value = 'Skinnlav'
    $ws.Range("G64").Value = 'Leptogium saturninum'
    $ws.Range("H64").Value = '(Dicks.) Nyl.'
    $ws.Range("Q64").Value = 838023
    $ws.Range("R64").Value = 7446298

    # Row 65
    $ws.Range("A65").Value = 112081041
    $ws.Range("B65").Value = 95693
    $ws.Range("E65").Value = 221941
    $ws.Range("F65").Value = 'Plattlummer'
    $ws.Range("G65").Value = 'Lycopodium complanatum'
    $ws.Range("H65").Value = 'L.'
    $ws.Range("Q65").Value = 837955
    $ws.Range("R65").Value = 7446342
    $ws.Range("Y65").Value = '2023-08-31'
    $ws.Range("AA65").Value = '2023-08-31'

    # Row 66
    $ws.Range("A66").Value = 112081046
    $ws.Range("B66").Value = 90792
    $ws.Range("E66").Value = 4361
    $ws.Range("F66").Value = 'Orange taggsvamp'
    $ws.Range("G66").Value = 'Hydnellum aurantiacum'
    $ws.Range("H66").Value = '(Batsch:Fr.) P.Karst.'
    $ws.Range("Q66").Value = 837922
    $ws.Range("R66").Value = 7446268
    $ws.Range("Y66").Value = '2023-08-31'
    $ws.Range("AA66").Value = '2023-08-31'

    # Row 67
    $ws.Range("A67").Value = 112080967
    $ws.Range("Q67").Value = 838146
    $ws.Range("R67").Value = 7446346
    $ws.Range("Y67").Value = '2023-09-06'
    $ws.Range("AA67").Value = '2023-09-06'

    # Row 68
    $ws.Range("A68").Value = 112081076
    $ws.Range("B68").Value = 90800
    $ws.Range("D68").Value = 'LC'
    $ws.Range("E68").Value = 4364
    $ws.Range("F68").Value = 'Dropptaggsvamp'
    $ws.Range("G68").Value = 'Hydnellum ferrugineum'
    $ws.Range("H68").Value = '(Fr.:Fr.) P. Karst.'
    $ws.Range("Q68").Value = 837900
    $ws.Range("R68").Value = 7446196
    $ws.Range("Y68").Value = '2023-09-06'
    $ws.Range("AA68").Value = '2023-09-06'

    # Row 69
    $ws.Range("A69").Value = 112081067
    $ws.Range("B69").Value = 77636
    $ws.Range("D69").Value = 'NT'
    $ws.Range("E69").Value = 6425
    $ws.Range("F69").Value = 'Garnlav'
    $ws.Range("G69").Value = 'Alectoria sarmentosa'
    $ws.Range("H69").Value = '(Ach.) Ach.'
    $ws.Range("Q69").Value = 838039
    $ws.Range("R69").Value = 7446399

    # Row 70
    $ws.Range("A70").Value = 112080980
    $ws.Range("B70").Value = 86357
    $ws.Range("D70").Value = 'NT'
    $ws.Range("E70").Value = 4412
    $ws.Range("F70").Value = 'Äggvaxskivling'
    $ws.Range("G70").Value = 'Hygrophorus karstenii'
    $ws.Range("H70").Value = 'Sacc. & Cub.'
    $ws.Range("Q70").Value = 837988
    $ws.Range("R70").Value = 7446497

    # Row 71
    $ws.Range("A71").Value = 112081085
    $ws.Range("B71").Value = 90794
    $ws.Range("E71").Value = 4362
    $ws.Range("F71").Value = 'Blå taggsvamp'
    $ws.Range("G71").Value = 'Hydnellum caeruleum'
    $ws.Range("H71").Value = '(Hornem.) P.Karst.'
    $ws.Range("Q71").Value = 837853
    $ws.Range("R71").Value = 7446288

